$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Most Dailies"
$ws.Range("A3").Value = "Hilts Trader Limited Time Items"
$ws.Range("A4").Value = "Hilts Trader Limited Time Items"
$ws.Range("A5").Value = "Weekly Reset"

$ws.Range("A5").Select()
